# Wheelchair BOM Revision 2012A
# Adds a D-sub connector plug line to the "Connector to front" group (row 17)
# and a new "D-sub standoffs" line item (row 34), pushing the grand-total
# SUM formula down to row 35.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Connector")

# --- Row 17: new "Plug" line for the D-sub connector (part of the merged
#     "Connector to front" category spanning A16:A19) ---
# Note: F17 is written before D17 so that new shared-string entries are
# created in supplier-part-number-first order (538-39-01-2240, then
# 39-01-2240), matching how the workbook was authored.
$ws.Range("C17").Value = "Molex"
$ws.Range("F17").Value = "538-39-01-2240"
$ws.Range("D17").Value = "39-01-2240"
$ws.Range("E17").Value = "Mouser"
$ws.Range("G17").Value = 1.25
$ws.Range("H17").Value = 1

# --- Row 34: new "D-sub standoffs" line item ---
$ws.Range("B34").Value = "D-sub standoffs"
$ws.Range("E34").Value = "McMaster-Carr"
$ws.Range("F34").Value = "93620A701"
$ws.Range("G34").Value = 1.4
$ws.Range("H34").Value = 4
$ws.Range("I34").Formula = "=G34*H34"

# --- Row 35: grand total moves down one row; formula text is re-entered
#     unchanged (still summing I2:I33, the original item range) ---
$ws.Range("I35").Formula = "=SUM(I2:I33)"

# --- Selection cursor, as left by the editor after finishing the edit ---
[void]$ws.Range("E23:E32").Select()
